# Populate Sheet1 with three rows of reference data (A:D), matching the
# uploaded "TEMP_XLS" sheet: cols A-C are numeric codes, col D is an
# "IMPE" text id.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 7
$ws.Range("B1").Value = 16
$ws.Range("C1").Value = 1101
$ws.Range("D1").Value = "IMPE1001"

$ws.Range("A2").Value = 8
$ws.Range("B2").Value = 16
$ws.Range("C2").Value = 1101
$ws.Range("D2").Value = "IMPE1002"

$ws.Range("A3").Value = 9
$ws.Range("B3").Value = 16
$ws.Range("C3").Value = 1101
$ws.Range("D3").Value = "IMPE1003"

# Leave the selection where the author's saved file shows it (one column
# past the data, row 1).
$ws.Range("E1").Select()
